# Adds the "Find" worksheet (Matlab/Python index comparison) as the last sheet
# and makes it the active sheet, matching the target workbook state.
$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Find"

# Move the new sheet to be the last tab (after SUM)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-fetch the worksheet reference by name: Move() reseats object handles
# by position, so the old $newSheet variable would now point at whatever
# sheet occupies its old slot (SUM) rather than the moved "Find" sheet.
$ws = $wb.Worksheets.Item("Find")

# Row 1: "Matlab" label + sequence 1..118 (skipping 7 and 112)
$ws.Cells.Item(1, 1).Value = "Matlab"
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 8
$ws.Cells.Item(1, 9).Value = 9
$ws.Cells.Item(1, 10).Value = 10
$ws.Cells.Item(1, 11).Value = 11
$ws.Cells.Item(1, 12).Value = 12
$ws.Cells.Item(1, 13).Value = 13
$ws.Cells.Item(1, 14).Value = 14
$ws.Cells.Item(1, 15).Value = 15
$ws.Cells.Item(1, 16).Value = 16
$ws.Cells.Item(1, 17).Value = 17
$ws.Cells.Item(1, 18).Value = 18
$ws.Cells.Item(1, 19).Value = 19
$ws.Cells.Item(1, 20).Value = 20
$ws.Cells.Item(1, 21).Value = 21
$ws.Cells.Item(1, 22).Value = 22
$ws.Cells.Item(1, 23).Value = 23
$ws.Cells.Item(1, 24).Value = 24
$ws.Cells.Item(1, 25).Value = 25
$ws.Cells.Item(1, 26).Value = 26
$ws.Cells.Item(1, 27).Value = 27
$ws.Cells.Item(1, 28).Value = 28
$ws.Cells.Item(1, 29).Value = 29
$ws.Cells.Item(1, 30).Value = 30
$ws.Cells.Item(1, 31).Value = 31
$ws.Cells.Item(1, 32).Value = 32
$ws.Cells.Item(1, 33).Value = 33
$ws.Cells.Item(1, 34).Value = 34
$ws.Cells.Item(1, 35).Value = 35
$ws.Cells.Item(1, 36).Value = 36
$ws.Cells.Item(1, 37).Value = 37
$ws.Cells.Item(1, 38).Value = 38
$ws.Cells.Item(1, 39).Value = 39
$ws.Cells.Item(1, 40).Value = 40
$ws.Cells.Item(1, 41).Value = 41
$ws.Cells.Item(1, 42).Value = 42
$ws.Cells.Item(1, 43).Value = 43
$ws.Cells.Item(1, 44).Value = 44
$ws.Cells.Item(1, 45).Value = 45
$ws.Cells.Item(1, 46).Value = 46
$ws.Cells.Item(1, 47).Value = 47
$ws.Cells.Item(1, 48).Value = 48
$ws.Cells.Item(1, 49).Value = 49
$ws.Cells.Item(1, 50).Value = 50
$ws.Cells.Item(1, 51).Value = 51
$ws.Cells.Item(1, 52).Value = 52
$ws.Cells.Item(1, 53).Value = 53
$ws.Cells.Item(1, 54).Value = 54
$ws.Cells.Item(1, 55).Value = 55
$ws.Cells.Item(1, 56).Value = 56
$ws.Cells.Item(1, 57).Value = 57
$ws.Cells.Item(1, 58).Value = 58
$ws.Cells.Item(1, 59).Value = 59
$ws.Cells.Item(1, 60).Value = 60
$ws.Cells.Item(1, 61).Value = 61
$ws.Cells.Item(1, 62).Value = 62
$ws.Cells.Item(1, 63).Value = 63
$ws.Cells.Item(1, 64).Value = 64
$ws.Cells.Item(1, 65).Value = 65
$ws.Cells.Item(1, 66).Value = 66
$ws.Cells.Item(1, 67).Value = 67
$ws.Cells.Item(1, 68).Value = 68
$ws.Cells.Item(1, 69).Value = 69
$ws.Cells.Item(1, 70).Value = 70
$ws.Cells.Item(1, 71).Value = 71
$ws.Cells.Item(1, 72).Value = 72
$ws.Cells.Item(1, 73).Value = 73
$ws.Cells.Item(1, 74).Value = 74
$ws.Cells.Item(1, 75).Value = 75
$ws.Cells.Item(1, 76).Value = 76
$ws.Cells.Item(1, 77).Value = 77
$ws.Cells.Item(1, 78).Value = 78
$ws.Cells.Item(1, 79).Value = 79
$ws.Cells.Item(1, 80).Value = 80
$ws.Cells.Item(1, 81).Value = 81
$ws.Cells.Item(1, 82).Value = 82
$ws.Cells.Item(1, 83).Value = 83
$ws.Cells.Item(1, 84).Value = 84
$ws.Cells.Item(1, 85).Value = 85
$ws.Cells.Item(1, 86).Value = 86
$ws.Cells.Item(1, 87).Value = 87
$ws.Cells.Item(1, 88).Value = 88
$ws.Cells.Item(1, 89).Value = 89
$ws.Cells.Item(1, 90).Value = 90
$ws.Cells.Item(1, 91).Value = 91
$ws.Cells.Item(1, 92).Value = 92
$ws.Cells.Item(1, 93).Value = 93
$ws.Cells.Item(1, 94).Value = 94
$ws.Cells.Item(1, 95).Value = 95
$ws.Cells.Item(1, 96).Value = 96
$ws.Cells.Item(1, 97).Value = 97
$ws.Cells.Item(1, 98).Value = 98
$ws.Cells.Item(1, 99).Value = 99
$ws.Cells.Item(1, 100).Value = 100
$ws.Cells.Item(1, 101).Value = 101
$ws.Cells.Item(1, 102).Value = 102
$ws.Cells.Item(1, 103).Value = 103
$ws.Cells.Item(1, 104).Value = 104
$ws.Cells.Item(1, 105).Value = 105
$ws.Cells.Item(1, 106).Value = 106
$ws.Cells.Item(1, 107).Value = 107
$ws.Cells.Item(1, 108).Value = 108
$ws.Cells.Item(1, 109).Value = 109
$ws.Cells.Item(1, 110).Value = 110
$ws.Cells.Item(1, 111).Value = 111
$ws.Cells.Item(1, 112).Value = 113
$ws.Cells.Item(1, 113).Value = 114
$ws.Cells.Item(1, 114).Value = 115
$ws.Cells.Item(1, 115).Value = 116
$ws.Cells.Item(1, 116).Value = 117
$ws.Cells.Item(1, 117).Value = 118

# Row 2: "Python" label + sequence 0..117 (skipping 6 and 111) = Row1 - 1
$ws.Cells.Item(2, 1).Value = "Python"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(2, 7).Value = 5
$ws.Cells.Item(2, 8).Value = 7
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(2, 12).Value = 11
$ws.Cells.Item(2, 13).Value = 12
$ws.Cells.Item(2, 14).Value = 13
$ws.Cells.Item(2, 15).Value = 14
$ws.Cells.Item(2, 16).Value = 15
$ws.Cells.Item(2, 17).Value = 16
$ws.Cells.Item(2, 18).Value = 17
$ws.Cells.Item(2, 19).Value = 18
$ws.Cells.Item(2, 20).Value = 19
$ws.Cells.Item(2, 21).Value = 20
$ws.Cells.Item(2, 22).Value = 21
$ws.Cells.Item(2, 23).Value = 22
$ws.Cells.Item(2, 24).Value = 23
$ws.Cells.Item(2, 25).Value = 24
$ws.Cells.Item(2, 26).Value = 25
$ws.Cells.Item(2, 27).Value = 26
$ws.Cells.Item(2, 28).Value = 27
$ws.Cells.Item(2, 29).Value = 28
$ws.Cells.Item(2, 30).Value = 29
$ws.Cells.Item(2, 31).Value = 30
$ws.Cells.Item(2, 32).Value = 31
$ws.Cells.Item(2, 33).Value = 32
$ws.Cells.Item(2, 34).Value = 33
$ws.Cells.Item(2, 35).Value = 34
$ws.Cells.Item(2, 36).Value = 35
$ws.Cells.Item(2, 37).Value = 36
$ws.Cells.Item(2, 38).Value = 37
$ws.Cells.Item(2, 39).Value = 38
$ws.Cells.Item(2, 40).Value = 39
$ws.Cells.Item(2, 41).Value = 40
$ws.Cells.Item(2, 42).Value = 41
$ws.Cells.Item(2, 43).Value = 42
$ws.Cells.Item(2, 44).Value = 43
$ws.Cells.Item(2, 45).Value = 44
$ws.Cells.Item(2, 46).Value = 45
$ws.Cells.Item(2, 47).Value = 46
$ws.Cells.Item(2, 48).Value = 47
$ws.Cells.Item(2, 49).Value = 48
$ws.Cells.Item(2, 50).Value = 49
$ws.Cells.Item(2, 51).Value = 50
$ws.Cells.Item(2, 52).Value = 51
$ws.Cells.Item(2, 53).Value = 52
$ws.Cells.Item(2, 54).Value = 53
$ws.Cells.Item(2, 55).Value = 54
$ws.Cells.Item(2, 56).Value = 55
$ws.Cells.Item(2, 57).Value = 56
$ws.Cells.Item(2, 58).Value = 57
$ws.Cells.Item(2, 59).Value = 58
$ws.Cells.Item(2, 60).Value = 59
$ws.Cells.Item(2, 61).Value = 60
$ws.Cells.Item(2, 62).Value = 61
$ws.Cells.Item(2, 63).Value = 62
$ws.Cells.Item(2, 64).Value = 63
$ws.Cells.Item(2, 65).Value = 64
$ws.Cells.Item(2, 66).Value = 65
$ws.Cells.Item(2, 67).Value = 66
$ws.Cells.Item(2, 68).Value = 67
$ws.Cells.Item(2, 69).Value = 68
$ws.Cells.Item(2, 70).Value = 69
$ws.Cells.Item(2, 71).Value = 70
$ws.Cells.Item(2, 72).Value = 71
$ws.Cells.Item(2, 73).Value = 72
$ws.Cells.Item(2, 74).Value = 73
$ws.Cells.Item(2, 75).Value = 74
$ws.Cells.Item(2, 76).Value = 75
$ws.Cells.Item(2, 77).Value = 76
$ws.Cells.Item(2, 78).Value = 77
$ws.Cells.Item(2, 79).Value = 78
$ws.Cells.Item(2, 80).Value = 79
$ws.Cells.Item(2, 81).Value = 80
$ws.Cells.Item(2, 82).Value = 81
$ws.Cells.Item(2, 83).Value = 82
$ws.Cells.Item(2, 84).Value = 83
$ws.Cells.Item(2, 85).Value = 84
$ws.Cells.Item(2, 86).Value = 85
$ws.Cells.Item(2, 87).Value = 86
$ws.Cells.Item(2, 88).Value = 87
$ws.Cells.Item(2, 89).Value = 88
$ws.Cells.Item(2, 90).Value = 89
$ws.Cells.Item(2, 91).Value = 90
$ws.Cells.Item(2, 92).Value = 91
$ws.Cells.Item(2, 93).Value = 92
$ws.Cells.Item(2, 94).Value = 93
$ws.Cells.Item(2, 95).Value = 94
$ws.Cells.Item(2, 96).Value = 95
$ws.Cells.Item(2, 97).Value = 96
$ws.Cells.Item(2, 98).Value = 97
$ws.Cells.Item(2, 99).Value = 98
$ws.Cells.Item(2, 100).Value = 99
$ws.Cells.Item(2, 101).Value = 100
$ws.Cells.Item(2, 102).Value = 101
$ws.Cells.Item(2, 103).Value = 102
$ws.Cells.Item(2, 104).Value = 103
$ws.Cells.Item(2, 105).Value = 104
$ws.Cells.Item(2, 106).Value = 105
$ws.Cells.Item(2, 107).Value = 106
$ws.Cells.Item(2, 108).Value = 107
$ws.Cells.Item(2, 109).Value = 108
$ws.Cells.Item(2, 110).Value = 109
$ws.Cells.Item(2, 111).Value = 110
$ws.Cells.Item(2, 112).Value = 112
$ws.Cells.Item(2, 113).Value = 113
$ws.Cells.Item(2, 114).Value = 114
$ws.Cells.Item(2, 115).Value = 115
$ws.Cells.Item(2, 116).Value = 116
$ws.Cells.Item(2, 117).Value = 117

# Row 3: BENER/SALAH check formula for each column B..DM
$ws.Cells.Item(3, 2).Formula = '=IF((B1)=(B2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 3).Formula = '=IF((C1)=(C2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 4).Formula = '=IF((D1)=(D2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 5).Formula = '=IF((E1)=(E2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 6).Formula = '=IF((F1)=(F2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 7).Formula = '=IF((G1)=(G2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 8).Formula = '=IF((H1)=(H2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 9).Formula = '=IF((I1)=(I2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 10).Formula = '=IF((J1)=(J2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 11).Formula = '=IF((K1)=(K2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 12).Formula = '=IF((L1)=(L2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 13).Formula = '=IF((M1)=(M2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 14).Formula = '=IF((N1)=(N2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 15).Formula = '=IF((O1)=(O2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 16).Formula = '=IF((P1)=(P2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 17).Formula = '=IF((Q1)=(Q2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 18).Formula = '=IF((R1)=(R2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 19).Formula = '=IF((S1)=(S2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 20).Formula = '=IF((T1)=(T2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 21).Formula = '=IF((U1)=(U2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 22).Formula = '=IF((V1)=(V2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 23).Formula = '=IF((W1)=(W2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 24).Formula = '=IF((X1)=(X2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 25).Formula = '=IF((Y1)=(Y2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 26).Formula = '=IF((Z1)=(Z2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 27).Formula = '=IF((AA1)=(AA2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 28).Formula = '=IF((AB1)=(AB2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 29).Formula = '=IF((AC1)=(AC2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 30).Formula = '=IF((AD1)=(AD2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 31).Formula = '=IF((AE1)=(AE2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 32).Formula = '=IF((AF1)=(AF2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 33).Formula = '=IF((AG1)=(AG2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 34).Formula = '=IF((AH1)=(AH2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 35).Formula = '=IF((AI1)=(AI2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 36).Formula = '=IF((AJ1)=(AJ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 37).Formula = '=IF((AK1)=(AK2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 38).Formula = '=IF((AL1)=(AL2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 39).Formula = '=IF((AM1)=(AM2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 40).Formula = '=IF((AN1)=(AN2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 41).Formula = '=IF((AO1)=(AO2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 42).Formula = '=IF((AP1)=(AP2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 43).Formula = '=IF((AQ1)=(AQ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 44).Formula = '=IF((AR1)=(AR2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 45).Formula = '=IF((AS1)=(AS2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 46).Formula = '=IF((AT1)=(AT2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 47).Formula = '=IF((AU1)=(AU2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 48).Formula = '=IF((AV1)=(AV2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 49).Formula = '=IF((AW1)=(AW2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 50).Formula = '=IF((AX1)=(AX2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 51).Formula = '=IF((AY1)=(AY2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 52).Formula = '=IF((AZ1)=(AZ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 53).Formula = '=IF((BA1)=(BA2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 54).Formula = '=IF((BB1)=(BB2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 55).Formula = '=IF((BC1)=(BC2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 56).Formula = '=IF((BD1)=(BD2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 57).Formula = '=IF((BE1)=(BE2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 58).Formula = '=IF((BF1)=(BF2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 59).Formula = '=IF((BG1)=(BG2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 60).Formula = '=IF((BH1)=(BH2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 61).Formula = '=IF((BI1)=(BI2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 62).Formula = '=IF((BJ1)=(BJ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 63).Formula = '=IF((BK1)=(BK2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 64).Formula = '=IF((BL1)=(BL2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 65).Formula = '=IF((BM1)=(BM2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 66).Formula = '=IF((BN1)=(BN2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 67).Formula = '=IF((BO1)=(BO2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 68).Formula = '=IF((BP1)=(BP2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 69).Formula = '=IF((BQ1)=(BQ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 70).Formula = '=IF((BR1)=(BR2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 71).Formula = '=IF((BS1)=(BS2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 72).Formula = '=IF((BT1)=(BT2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 73).Formula = '=IF((BU1)=(BU2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 74).Formula = '=IF((BV1)=(BV2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 75).Formula = '=IF((BW1)=(BW2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 76).Formula = '=IF((BX1)=(BX2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 77).Formula = '=IF((BY1)=(BY2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 78).Formula = '=IF((BZ1)=(BZ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 79).Formula = '=IF((CA1)=(CA2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 80).Formula = '=IF((CB1)=(CB2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 81).Formula = '=IF((CC1)=(CC2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 82).Formula = '=IF((CD1)=(CD2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 83).Formula = '=IF((CE1)=(CE2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 84).Formula = '=IF((CF1)=(CF2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 85).Formula = '=IF((CG1)=(CG2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 86).Formula = '=IF((CH1)=(CH2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 87).Formula = '=IF((CI1)=(CI2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 88).Formula = '=IF((CJ1)=(CJ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 89).Formula = '=IF((CK1)=(CK2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 90).Formula = '=IF((CL1)=(CL2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 91).Formula = '=IF((CM1)=(CM2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 92).Formula = '=IF((CN1)=(CN2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 93).Formula = '=IF((CO1)=(CO2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 94).Formula = '=IF((CP1)=(CP2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 95).Formula = '=IF((CQ1)=(CQ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 96).Formula = '=IF((CR1)=(CR2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 97).Formula = '=IF((CS1)=(CS2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 98).Formula = '=IF((CT1)=(CT2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 99).Formula = '=IF((CU1)=(CU2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 100).Formula = '=IF((CV1)=(CV2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 101).Formula = '=IF((CW1)=(CW2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 102).Formula = '=IF((CX1)=(CX2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 103).Formula = '=IF((CY1)=(CY2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 104).Formula = '=IF((CZ1)=(CZ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 105).Formula = '=IF((DA1)=(DA2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 106).Formula = '=IF((DB1)=(DB2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 107).Formula = '=IF((DC1)=(DC2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 108).Formula = '=IF((DD1)=(DD2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 109).Formula = '=IF((DE1)=(DE2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 110).Formula = '=IF((DF1)=(DF2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 111).Formula = '=IF((DG1)=(DG2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 112).Formula = '=IF((DH1)=(DH2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 113).Formula = '=IF((DI1)=(DI2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 114).Formula = '=IF((DJ1)=(DJ2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 115).Formula = '=IF((DK1)=(DK2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 116).Formula = '=IF((DL1)=(DL2+1), "BENER", "SALAH")'
$ws.Cells.Item(3, 117).Formula = '=IF((DM1)=(DM2+1), "BENER", "SALAH")'

# Select the cell that was active in the target workbook, then activate the sheet
$ws.Range("DG3").Select()
$ws.Activate()
